$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A12 with new text values (all as text / strings, not numbers)
$ws.Range("A2").Value = "E0000000"
$ws.Range("A3").Value = "E1111111"
$ws.Range("A4").Value = "E1212121"
$ws.Range("A5").Value = "40000000"
$ws.Range("A6").Value = "41111111"
$ws.Range("A7").Value = "41212121"
$ws.Range("A8").Value = "E6942000"
$ws.Range("A9").Value = "43333333"
$ws.Range("A10").Value = "E8787655"
$ws.Range("A11").Value = "E3213214"
$ws.Range("A12").Value = "49856473"

# Update the selection state to match diff (E6:E7, active cell E7)
$ws.Range("E6:E7").Select()
